# "Ajout d'une date de visite #51"
#
# The "tableauBatiment" table (shape 1) on slide 2 gets a new row
# "Date de visite :" / "Date de visite" appended after "Date de
# rénovation", matching the formatting of the preceding rows. The table
# frame grows to fit the new row, and the context/title bar shape below
# the photo ("Elements de contexte sur le bâtiment", shape 5) is pushed
# down to make room for it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$tableShape = $s.Shapes.Item(1)
$tbl = $tableShape.Table

# Append a new row; PowerPoint seeds it by copying the formatting of the
# last existing row ("Date de rénovation").
$null = $tbl.Rows.Add()
$lastRow = $tbl.Rows.Count

$tbl.Cell($lastRow, 1).Shape.TextFrame.TextRange.Text = "Date de visite :"
$tbl.Cell($lastRow, 2).Shape.TextFrame.TextRange.Text = "Date de visite"

# Resync the graphic frame's stored height with the table's new
# (autofit) rendered height now that it has six rows.
$tableShape.Height = 231.5832

# Shift the title/context bar shape down below the taller table.
$contextShape = $s.Shapes.Item(5)
$contextShape.Left = 544.44336
$contextShape.Top = 373.63414
